$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1469.7142
$ws.Range("I2").Value = 337.66666
$ws.Range("J2").Value = 2318.75
$ws.Range("K2").Value = 337.66666
$ws.Range("L2").Value = 2318.75
$ws.Range("M2").Value = -224.66666
$ws.Range("N2").Value = -2544.75
$ws.Range("H98").Value = 2135.2
$ws.Range("I98").Value = 2135.2
$ws.Range("K98").Value = 2135.2
$ws.Range("M98").Value = -637.1999999999998
$ws.Range("H100").Value = 1905.2333
$ws.Range("I100").Value = 1411
$ws.Range("K100").Value = 1411
$ws.Range("M100").Value = -870
$ws.Range("H112").Value = 2633.1614
$ws.Range("J112").Value = 2643.541
$ws.Range("L112").Value = 7930.623000000001
$ws.Range("N112").Value = -10146.623
$ws.Range("H122").Value = 2135.2
$ws.Range("I122").Value = 2135.2
$ws.Range("K122").Value = 6405.599999999999
$ws.Range("M122").Value = -3955.599999999999
$ws.Range("H129").Value = 4135.125
$ws.Range("I129").Value = 627.2857
$ws.Range("K129").Value = 1881.8571
$ws.Range("M129").Value = 3118.1429
$ws.Range("H138").Value = 5071.923
$ws.Range("I138").Value = 1788.258
$ws.Range("K138").Value = 5364.774
$ws.Range("M138").Value = -224.7740000000003

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2563.8333
$ws.Range("J45").Value = 4405.2856
$ws.Range("L45").Value = 4405.2856
$ws.Range("N45").Value = -5159.2856
$ws.Range("H64").Value = 112499.5
$ws.Range("J64").Value = 112499.5
$ws.Range("L64").Value = 112499.5
$ws.Range("N64").Value = -112995.5
$ws.Range("H67").Value = 112499.5
$ws.Range("J67").Value = 112499.5
$ws.Range("L67").Value = 112499.5
$ws.Range("N67").Value = -114215.5
$ws.Range("H88").Value = 85799.664
$ws.Range("I88").Value = 1199.6666
$ws.Range("K88").Value = 1199.6666
$ws.Range("M88").Value = -793.6666
$ws.Range("H91").Value = 85799.664
$ws.Range("I91").Value = 1199.6666
$ws.Range("K91").Value = 1199.6666
$ws.Range("M91").Value = 204.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10869.786
$ws.Range("I20").Value = 12247.375
$ws.Range("K20").Value = 12247.375
$ws.Range("M20").Value = -12000.375
$ws.Range("H105").Value = 2701.5454
$ws.Range("I105").Value = 2635.2222
$ws.Range("K105").Value = 2635.2222
$ws.Range("M105").Value = -888.2222000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 638.9167
$ws.Range("I107").Value = 496.33334
$ws.Range("J107").Value = 1066.6666
$ws.Range("K107").Value = 496.33334
$ws.Range("L107").Value = 1066.6666
$ws.Range("M107").Value = 1423.66666
$ws.Range("N107").Value = -4906.6666
$ws.Range("H132").Value = 27788268
$ws.Range("I132").Value = 37044684
$ws.Range("K132").Value = 111134052
$ws.Range("M132").Value = -111131522
$ws.Range("H140").Value = 105296.664
$ws.Range("J140").Value = 105296.664
$ws.Range("L140").Value = 105296.664
$ws.Range("N140").Value = -115656.664

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27207882
$ws.Range("I4").Value = 653620.9399999999
$ws.Range("J4").Value = 76522936
$ws.Range("K4").Value = 1960862.82
$ws.Range("L4").Value = 229568808
$ws.Range("M4").Value = -1960750.82
$ws.Range("N4").Value = -229569032
$ws.Range("H131").Value = 19667614
$ws.Range("I131").Value = 27778580
$ws.Range("J131").Value = 18480642
$ws.Range("K131").Value = 83335740
$ws.Range("L131").Value = 55441926
$ws.Range("M131").Value = -83330700
$ws.Range("N131").Value = -55452006
$ws.Range("H140").Value = 8503.9375
$ws.Range("I140").Value = 2049.4348
$ws.Range("J140").Value = 24998.777
$ws.Range("K140").Value = 6148.3044
$ws.Range("L140").Value = 74996.33099999999
$ws.Range("M140").Value = -968.3044
$ws.Range("N140").Value = -85356.33099999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8465.444
$ws.Range("I70").Value = 9073
$ws.Range("K70").Value = 9073
$ws.Range("M70").Value = -8803
$ws.Range("H73").Value = 8465.444
$ws.Range("I73").Value = 9073
$ws.Range("K73").Value = 9073
$ws.Range("M73").Value = -8137
$ws.Range("H122").Value = 483008.12
$ws.Range("I122").Value = 918891.75
$ws.Range("J122").Value = 7498.727
$ws.Range("K122").Value = 2756675.25
$ws.Range("L122").Value = 22496.181
$ws.Range("M122").Value = -2754225.25
$ws.Range("N122").Value = -27396.181
$ws.Range("H126").Value = 3681.7568
$ws.Range("I126").Value = 2538.111
$ws.Range("J126").Value = 4765.2104
$ws.Range("K126").Value = 7614.333
$ws.Range("L126").Value = 14295.6312
$ws.Range("M126").Value = -5144.333
$ws.Range("N126").Value = -19235.6312

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4249.683
$ws.Range("I7").Value = 3096.35
$ws.Range("J7").Value = 5348.095
$ws.Range("K7").Value = 3096.35
$ws.Range("L7").Value = 5348.095
$ws.Range("M7").Value = -2984.35
$ws.Range("N7").Value = -5572.095
$ws.Range("H40").Value = 4313833.5
$ws.Range("I40").Value = 6251503
$ws.Range("K40").Value = 6251503
$ws.Range("M40").Value = -6251367
$ws.Range("H46").Value = 3541.3333
$ws.Range("I46").Value = 926
$ws.Range("K46").Value = 926
$ws.Range("M46").Value = -738
$ws.Range("H122").Value = 6535.5264
$ws.Range("I122").Value = 4422.591
$ws.Range("J122").Value = 9440.8125
$ws.Range("K122").Value = 13267.773
$ws.Range("L122").Value = 28322.4375
$ws.Range("M122").Value = -10817.773
$ws.Range("N122").Value = -33222.4375
$ws.Range("H126").Value = 4249.683
$ws.Range("I126").Value = 3096.35
$ws.Range("J126").Value = 5348.095
$ws.Range("K126").Value = 9289.049999999999
$ws.Range("L126").Value = 16044.285
$ws.Range("M126").Value = -6819.049999999999
$ws.Range("N126").Value = -20984.285
$ws.Range("H136").Value = 4037.8635
$ws.Range("I136").Value = 2345.5
$ws.Range("K136").Value = 7036.5
$ws.Range("M136").Value = -4486.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2454.923
$ws.Range("I100").Value = 2373.6333
$ws.Range("K100").Value = 4747.2666
$ws.Range("M100").Value = -4206.2666
$ws.Range("H122").Value = 5841.1035
$ws.Range("I122").Value = 5067.5557
$ws.Range("J122").Value = 7106.909
$ws.Range("K122").Value = 15202.6671
$ws.Range("L122").Value = 21320.727
$ws.Range("M122").Value = -12752.6671
$ws.Range("N122").Value = -26220.727
$ws.Range("H132").Value = 3006.491
$ws.Range("I132").Value = 1360.6389
$ws.Range("J132").Value = 6124.9473
$ws.Range("K132").Value = 4081.9167
$ws.Range("L132").Value = 18374.8419
$ws.Range("M132").Value = -1551.9167
$ws.Range("N132").Value = -23434.8419
$ws.Range("H136").Value = 7988.268
$ws.Range("I136").Value = 1588.6875
$ws.Range("K136").Value = 4766.0625
$ws.Range("M136").Value = -2216.0625
